# "Completed squaring off futures analysis"
#
# The "Hedging & Payoff" sheet (sheet index 1) had its existing content
# (originally rows 33-50) shifted up by 25 rows (to rows 8-25) by deleting
# rows 1-25, a new section header was added at the (now empty) top of the
# sheet, and a brand-new "Squaring off Futures" walkthrough + payoff
# calculation block was appended below the existing payoff table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hedging & Payoff")

# --- Step 1: shift the existing content from rows 33-50 up to rows 8-25 ---
# by deleting the (empty) rows 1-25 above it.
$ws.Range("A1:A25").EntireRow.Delete()

# --- Step 2: add the new section title above the existing content ---
$ws.Range("B7").Value = "Trading Data and Payoff"
$ws.Range("B7").Font.Bold = $true

# --- Step 3: build the new "Squaring off Futures" section below row 25 ---
$ws.Range("B28").Value = "Squaring off Futures"
$ws.Range("B28").Font.Bold = $true

$ws.Range("B29").Value = "Existing Futures contract without settlement (before expiry)"

$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "You expected the prices to increase, but they are decreasing"

$ws.Range("A31").Value = 2
$ws.Range("B31").Value = "or You no longer hold the asset"

$ws.Range("B34").Value = "thus you go Short Futures at price"
$ws.Range("B34").Font.Bold = $true
$ws.Range("C16").Copy()
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C34").Formula = "='Futures Data'!C41"
$ws.Range("D34").Formula = "=D16"

$ws.Range("B33").Value = "On 29th May, you want to exit"

$ws.Range("B36").Value = "29 May Pay off from Futures Profit"
$ws.Range("B36").Font.Bold = $true
$ws.Range("C36").Formula = "=C34-C16"
$ws.Range("C36").NumberFormat = "#,##0.00"
$ws.Range("D36").Formula = "=D34"

# --- Step 4: restore view state ---
# "Futures Data" loses its tab-selected flag and ends up scrolled down with
# B36 selected; "Hedging & Payoff" becomes the active tab with C36 selected.
$wsFutures = $wb.Worksheets.Item("Futures Data")
$wsFutures.Activate()
$wsFutures.Range("B36").Select()

$ws.Activate()
$ws.Range("C36").Select()
